# "Generate Report for Handoff"
# Updates the localization-status workbook: a new source markdown file
# (15cffeb1-...) replaces the old one (4beece0e-...) and two new PNG
# assets are introduced for handoff, each getting a fresh handoff
# package/time stamp. Applies to the Overview sheet and to both the
# zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$sourceRepo = "https://github.com/OpenLocalizationTest/oltest/blob/c7fab869f5b23efbe95e3dc85df59ccd84df8473"
$zhHandoffRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1256bfceff5776d1541734ad392bb3d084795a12/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7521a4735c74768a30edd2f102f822bc50df5153/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$mdGuid = "15cffeb1-2a3b-43fb-b7b8-dc385f233fd4"
$png1Guid = "41823f18-b08f-4f12-9095-20bca3177b4d"
$png2Guid = "e5701b81-15b7-479a-963d-bd0c7ea05c69"
$zhPngHash = "550b1b8a9ae01e4be2390dd543624bdecf393374"
$dePngHash = "a7484058ea4f5f4545c074de8f68e115e37413b4"
$xlfHash = "7b59110197c65b1853cc1afaa3eab2cf8d373c9d"

$mdDisplay = "$mdGuid.md"
$png1Display = "$png1Guid.png"
$png2Display = "$png2Guid.png"
$configDisplay = ".localization-config"

$readyStatus = "Ready for handoff"
$notLocalizedStatus = "Not to be localized"
$epochTime = "0001-01-01 00:00:00"
$dependencyFrom = "e2e\$mdGuid.md"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = $readyStatus
$ws1.Range("C2").Value = $readyStatus

$ws1.Range("B3").Value = $readyStatus
$ws1.Range("C3").Value = $readyStatus

$ws1.Range("B4").Value = $readyStatus
$ws1.Range("C4").Value = $readyStatus

$ws1.Range("B5").Value = $notLocalizedStatus
$ws1.Range("C5").Value = $notLocalizedStatus

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$sourceRepo/e2e/$mdDisplay", [Type]::Missing, [Type]::Missing, $mdDisplay)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$sourceRepo/e2e/$png1Display", [Type]::Missing, [Type]::Missing, $png1Display)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$sourceRepo/e2e/$png2Display", [Type]::Missing, [Type]::Missing, $png2Display)
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$sourceRepo/$configDisplay", [Type]::Missing, [Type]::Missing, $configDisplay)

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("B2").Value = $readyStatus
$ws2.Range("D2").Value = "2016-03-09 18:55:38"
$ws2.Range("G2").Value = $epochTime
$ws2.Range("H2").Value = "Include"

$ws2.Range("B3").Value = $readyStatus
$ws2.Range("D3").Value = "2016-03-09 18:55:38"
$ws2.Range("G3").Value = $epochTime
$ws2.Range("H3").Value = "IsDependency"
$ws2.Range("I3").Value = $dependencyFrom

$ws2.Range("B4").Value = $readyStatus
$ws2.Range("D4").Value = "2016-03-09 18:55:38"
$ws2.Range("G4").Value = $epochTime
$ws2.Range("H4").Value = "IsDependency"
$ws2.Range("I4").Value = $dependencyFrom

$ws2.Range("B5").Value = $notLocalizedStatus
$ws2.Range("D5").Value = $epochTime
$ws2.Range("G5").Value = $epochTime
$ws2.Range("H5").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$sourceRepo/e2e/$mdDisplay", [Type]::Missing, [Type]::Missing, $mdDisplay)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhHandoffRepo/$mdGuid.$xlfHash.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$mdGuid.$xlfHash.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$sourceRepo/e2e/$png1Display", [Type]::Missing, [Type]::Missing, $png1Display)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhHandoffRepo/$zhPngHash.png", [Type]::Missing, [Type]::Missing, "$zhPngHash.png")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$sourceRepo/e2e/$png2Display", [Type]::Missing, [Type]::Missing, $png2Display)
$ws2.Hyperlinks.Add($ws2.Range("C4"), "$zhHandoffRepo/$dePngHash.png", [Type]::Missing, [Type]::Missing, "$dePngHash.png")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "$sourceRepo/$configDisplay", [Type]::Missing, [Type]::Missing, $configDisplay)

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("B2").Value = $readyStatus
$ws3.Range("D2").Value = "2016-03-09 18:55:42"
$ws3.Range("G2").Value = $epochTime
$ws3.Range("H2").Value = "Include"

$ws3.Range("B3").Value = $readyStatus
$ws3.Range("D3").Value = "2016-03-09 18:55:42"
$ws3.Range("G3").Value = $epochTime
$ws3.Range("H3").Value = "IsDependency"
$ws3.Range("I3").Value = $dependencyFrom

$ws3.Range("B4").Value = $readyStatus
$ws3.Range("D4").Value = "2016-03-09 18:55:42"
$ws3.Range("G4").Value = $epochTime
$ws3.Range("H4").Value = "IsDependency"
$ws3.Range("I4").Value = $dependencyFrom

$ws3.Range("B5").Value = $notLocalizedStatus
$ws3.Range("D5").Value = $epochTime
$ws3.Range("G5").Value = $epochTime
$ws3.Range("H5").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$sourceRepo/e2e/$mdDisplay", [Type]::Missing, [Type]::Missing, $mdDisplay)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deHandoffRepo/$mdGuid.$xlfHash.de-de.xlf", [Type]::Missing, [Type]::Missing, "$mdGuid.$xlfHash.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$sourceRepo/e2e/$png1Display", [Type]::Missing, [Type]::Missing, $png1Display)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deHandoffRepo/$zhPngHash.png", [Type]::Missing, [Type]::Missing, "$zhPngHash.png")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$sourceRepo/e2e/$png2Display", [Type]::Missing, [Type]::Missing, $png2Display)
$ws3.Hyperlinks.Add($ws3.Range("C4"), "$deHandoffRepo/$dePngHash.png", [Type]::Missing, [Type]::Missing, "$dePngHash.png")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "$sourceRepo/$configDisplay", [Type]::Missing, [Type]::Missing, $configDisplay)

Write-Host "Report generated for handoff."
